$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7:A12").EntireRow.AutoFit()
$ws.Range("I6:L12").ClearContents()
